$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.90628671274171
$ws.Range("C2").Value = 5.57907341046594
$ws.Range("E2").Value = 9.400280221614485
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.761747239599984
$ws.Range("I2").Value = 40.24855956482254
$ws.Range("K2").Value = 14.55078682272254
$ws.Range("L2").Value = 10.63593255002355
$ws.Range("B3").Value = 15.84328852390039
$ws.Range("C3").Value = 5.386771106319451
$ws.Range("E3").Value = 9.414737312516177
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.764976007618043
$ws.Range("I3").Value = 39.87451427733728
$ws.Range("K3").Value = 14.49746023801181
$ws.Range("L3").Value = 10.6281118936766
$ws.Range("B4").Value = 15.80983635688796
$ws.Range("C4").Value = 5.266436261966263
$ws.Range("E4").Value = 9.425189516038218
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.767059842517659
$ws.Range("I4").Value = 39.6452182689921
$ws.Range("K4").Value = 14.46880992292949
$ws.Range("L4").Value = 10.62525600788562
$ws.Range("B5").Value = 15.79753187295052
$ws.Range("C5").Value = 5.216908484901765
$ws.Range("E5").Value = 9.429845136378617
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.767934604275964
$ws.Range("I5").Value = 39.5519229702403
$ws.Range("K5").Value = 14.45817310154288
$ws.Range("L5").Value = 10.62458259144507
$ws.Range("B6").Value = 15.7955692350668
$ws.Range("C6").Value = 5.208657226698295
$ws.Range("E6").Value = 9.430642133385355
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.768081405679852
$ws.Range("I6").Value = 39.53644158840267
$ws.Range("K6").Value = 14.45646982912805
$ws.Range("L6").Value = 10.62450041390904
$ws.Range("B7").Value = 15.80966502390055
$ws.Range("C7").Value = 5.265770190550938
$ws.Range("E7").Value = 9.425250698829075
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.767071536155242
$ws.Range("I7").Value = 39.64395940261377
$ws.Range("K7").Value = 14.46866225525817
$ws.Range("L7").Value = 10.62524493926019
$ws.Range("B8").Value = 15.88348786457602
$ws.Range("C8").Value = 5.513281554029345
$ws.Range("E8").Value = 9.404938105333978
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.762839539068709
$ws.Range("I8").Value = 40.1195434711484
$ws.Range("K8").Value = 14.5315564507346
$ws.Range("L8").Value = 10.63283278680075
$ws.Range("B9").Value = 16.06909434459523
$ws.Range("C9").Value = 5.977617918691116
$ws.Range("E9").Value = 9.377600397148942
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.755340509134915
$ws.Range("I9").Value = 41.05278463250566
$ws.Range("K9").Value = 14.68688992500535
$ws.Range("L9").Value = 10.66309794146423
$ws.Range("B10").Value = 16.22935880452259
$ws.Range("C10").Value = 6.302260254761609
$ws.Range("E10").Value = 9.365122940448346
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.750312546282561
$ws.Range("I10").Value = 41.73569650533796
$ws.Range("K10").Value = 14.81979854555842
$ws.Range("L10").Value = 10.69462416785826
$ws.Range("B11").Value = 16.30720782239771
$ws.Range("C11").Value = 6.44571373569083
$ws.Range("E11").Value = 9.361095119513339
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.748128468831503
$ws.Range("I11").Value = 42.04510568065547
$ws.Range("K11").Value = 14.8841615708235
$ws.Range("L11").Value = 10.71095901865309
$ws.Range("B12").Value = 16.33737485554462
$ws.Range("C12").Value = 6.49938167459814
$ws.Range("E12").Value = 9.359806475029933
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.747316151728793
$ws.Range("I12").Value = 42.16203421418483
$ws.Range("K12").Value = 14.90907844797059
$ws.Range("L12").Value = 10.71742867276784
$ws.Range("B13").Value = 16.33084766818485
$ws.Range("C13").Value = 6.487853131830337
$ws.Range("E13").Value = 9.360073492005514
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.747490444361758
$ws.Range("I13").Value = 42.13686306941513
$ws.Range("K13").Value = 14.90368824190364
$ws.Range("L13").Value = 10.71602273111875
$ws.Range("B14").Value = 16.30967604810609
$ws.Range("C14").Value = 6.450142400528443
$ws.Range("E14").Value = 9.360984362683416
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.748061344003401
$ws.Range("I14").Value = 42.05473056317379
$ws.Range("K14").Value = 14.88620069657886
$ws.Range("L14").Value = 10.71148560535267
$ws.Range("B15").Value = 16.29679660127366
$ws.Range("C15").Value = 6.426956922554885
$ws.Range("E15").Value = 9.361573096602507
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.748412954175266
$ws.Range("I15").Value = 42.00438931684424
$ws.Range("K15").Value = 14.87555938847082
$ws.Range("L15").Value = 10.70874338672638
$ws.Range("B16").Value = 16.22436868452836
$ws.Range("C16").Value = 6.292796026139897
$ws.Range("E16").Value = 9.36541930139756
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.750457349198439
$ws.Range("I16").Value = 41.7154457408009
$ws.Range("K16").Value = 14.81566935172631
$ws.Range("L16").Value = 10.6935965393285
$ws.Range("B17").Value = 16.18118639829079
$ws.Range("C17").Value = 6.209374687280295
$ws.Range("E17").Value = 9.368200708391724
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 3.75173787910685
$ws.Range("I17").Value = 41.53782638535344
$ws.Range("K17").Value = 14.77991634458027
$ws.Range("L17").Value = 10.68481320276539
$ws.Range("B18").Value = 16.15681594258299
$ws.Range("C18").Value = 6.160997413227346
$ws.Range("E18").Value = 9.369955669350972
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 3.752484122054359
$ws.Range("I18").Value = 41.43555015876117
$ws.Range("K18").Value = 14.7597205512252
$ws.Range("L18").Value = 10.67994906564673
$ws.Range("B19").Value = 16.14864540419912
$ws.Range("C19").Value = 6.14455128584455
$ws.Range("E19").Value = 9.37057653069434
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 3.752738458570976
$ws.Range("I19").Value = 41.40090333588621
$ws.Range("K19").Value = 14.75294637924523
$ws.Range("L19").Value = 10.67833448438659
$ws.Range("B20").Value = 16.18573507568108
$ws.Range("C20").Value = 6.218296323442988
$ws.Range("E20").Value = 9.367888566188633
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("G20").Value = 3.751600559596782
$ws.Range("I20").Value = 41.55674651787275
$ws.Range("K20").Value = 14.78368430716847
$ws.Range("L20").Value = 10.68572878508814
$ws.Range("B21").Value = 16.31587620815546
$ws.Range("C21").Value = 6.461237072574042
$ws.Range("E21").Value = 9.360710400173888
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.747893257461301
$ws.Range("I21").Value = 42.07886179208201
$ws.Range("K21").Value = 14.89132259025092
$ws.Range("L21").Value = 10.71281058396129
$ws.Range("B22").Value = 16.40492516992808
$ws.Range("C22").Value = 6.6161784436833
$ws.Range("E22").Value = 9.357397917805462
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.745556228615867
$ws.Range("I22").Value = 42.4186803629344
$ws.Range("K22").Value = 14.96483261377907
$ws.Range("L22").Value = 10.7321640553817
$ws.Range("B23").Value = 16.35704074742867
$ws.Range("C23").Value = 6.533848221479742
$ws.Range("E23").Value = 9.359039841179024
$ws.Range("F23").Value = 21.82633154475857
$ws.Range("G23").Value = 3.746795713927697
$ws.Range("I23").Value = 42.23746127293469
$ws.Range("K23").Value = 14.92531550985955
$ws.Range("L23").Value = 10.72168435326223
$ws.Range("B24").Value = 16.1836771976545
$ws.Range("C24").Value = 6.214264147741597
$ws.Range("E24").Value = 9.368029200071472
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.751662610437308
$ws.Range("I24").Value = 41.54819322138109
$ws.Range("K24").Value = 14.78197969147306
$ws.Range("L24").Value = 10.68531427225942
$ws.Range("B25").Value = 16.01460658609288
$ws.Range("C25").Value = 5.854659798170069
$ws.Range("E25").Value = 9.383659015600575
$ws.Range("F25").Value = 18.34778573295691
$ws.Range("G25").Value = 3.757284189093042
$ws.Range("I25").Value = 40.8005717363847
$ws.Range("K25").Value = 14.64151196931344
$ws.Range("L25").Value = 10.65327137114973
